$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'247.10"
$ws.Range("G2").Value = "'17"

# Row 3
$ws.Range("D3").Value = "'26.52"
$ws.Range("G3").Value = "'17"

# Row 4
$ws.Range("D4").Value = "'5.078"
$ws.Range("G4").Value = "'17"

# Row 5
$ws.Range("G5").Value = "'17"

# Row 6
$ws.Range("D6").Value = "'6.484"
$ws.Range("G6").Value = "'17"

# Row 7
$ws.Range("D7").Value = "'0.8131"
$ws.Range("G7").Value = "'17"

# Row 8
$ws.Range("D8").Value = "'0.8445"
$ws.Range("G8").Value = "'17"

# Row 9
$ws.Range("B9").Value = 'BitrueCoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D9").Value = "'0.02866"
$ws.Range("E9").Value = '8BitrueCoinBTR'
$ws.Range("G9").Value = "'17"

# Row 10
$ws.Range("B10").Value = 'BitMartToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D10").Value = "'0.09407"
$ws.Range("E10").Value = '9BitMartTokenBMX'
$ws.Range("G10").Value = "'17"

# Row 11
$ws.Range("B11").Value = 'BitForexToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D11").Value = "'0.001510"
$ws.Range("E11").Value = '10BitForexTokenBF'
$ws.Range("G11").Value = "'17"

# Row 12
$ws.Range("B12").Value = 'One'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D12").Value = "'0.0005999"
$ws.Range("E12").Value = '11OneONE'
$ws.Range("G12").Value = "'17"

# Row 13
$ws.Range("D13").Value = "'0.006191"
$ws.Range("G13").Value = "'17"

# Row 14
$ws.Range("D14").Value = "'3.594"
$ws.Range("G14").Value = "'17"

# Row 15
$ws.Range("D15").Value = "'3.011"
$ws.Range("G15").Value = "'17"

# Row 16
$ws.Range("G16").Value = "'17"

# Row 17
$ws.Range("G17").Value = "'17"

# Row 18
$ws.Range("G18").Value = "'17"

# Row 19
$ws.Range("D19").Value = "'0.06974"
$ws.Range("G19").Value = "'17"

# Row 20
$ws.Range("D20").Value = "'0.03202"
$ws.Range("G20").Value = "'17"

# Row 21
$ws.Range("G21").Value = "'17"

# Row 22
$ws.Range("D22").Value = "'3.739"
$ws.Range("G22").Value = "'17"

# Row 23
$ws.Range("D23").Value = "'0.04672"
$ws.Range("G23").Value = "'17"

# Row 24
$ws.Range("D24").Value = "'0.1350"
$ws.Range("G24").Value = "'17"

# Row 25
$ws.Range("D25").Value = "'0.001252"
$ws.Range("G25").Value = "'17"

# Row 26
$ws.Range("D26").Value = "'0.004601"
$ws.Range("G26").Value = "'17"

# Row 27
$ws.Range("D27").Value = "'0.00009600"
$ws.Range("G27").Value = "'17"

# Row 28
$ws.Range("G28").Value = "'17"

# Row 29
$ws.Range("G29").Value = "'17"

# Row 30
$ws.Range("G30").Value = "'17"

# Row 31
$ws.Range("G31").Value = "'17"

# Row 32
$ws.Range("G32").Value = "'17"

# Row 33
$ws.Range("G33").Value = "'17"

# Row 34
$ws.Range("G34").Value = "'17"

# Row 35
$ws.Range("G35").Value = "'17"

# Row 36
$ws.Range("G36").Value = "'17"

# Row 37
$ws.Range("G37").Value = "'17"

# Row 38
$ws.Range("G38").Value = "'17"

# Row 39
$ws.Range("G39").Value = "'17"

# Row 40
$ws.Range("D40").Value = "'0.03682"
$ws.Range("G40").Value = "'17"

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = "'0.1368"
$ws.Range("E41").Value = '40BKEXTokenBKKBestin24h'
$ws.Range("G41").Value = "'17"

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D42").Value = "'0.002660"
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("G42").Value = "'17"

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D43").Value = "'0.006220"
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("G43").Value = "'17"

# Row 44
$ws.Range("D44").Value = "'0.008912"
$ws.Range("G44").Value = "'17"

# Row 45
$ws.Range("D45").Value = "'0.00005294"
$ws.Range("G45").Value = "'17"

# Row 46
$ws.Range("G46").Value = "'17"

# Row 47
$ws.Range("D47").Value = "'0.1200"
$ws.Range("G47").Value = "'17"

# Row 48
$ws.Range("D48").Value = "'0.002521"
$ws.Range("G48").Value = "'17"

# Row 49
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("G49").Value = "'17"

# Row 50
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("G50").Value = "'17"

# Row 51
$ws.Range("G51").Value = "'17"
